$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.527.59"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.596.02"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.75"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.38"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.250"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0597"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.826.89"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.601.87"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.548.64"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("E16").Value = "  +3.32%  "
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.63"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.14"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.27%  "
$ws.Range("E20").Value = "  +2.83%  "
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.82"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.32"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.108"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("E31").Value = "  +3.17%  "
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.13"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.431.19"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.55"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.78%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.02"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.81"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.01%  "
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("E41").Value = "  +3.46%  "
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("E43").Value = "  +7.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "53.25"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +25.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.799"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.40%  "
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.985"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +18.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.53"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.17%  "
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.737.12"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "86.17"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.46%  "
